$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.624.05'
$ws.Range("D3").Value = '3.601.96'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.21%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '8.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D12").Value = '4.213.32'
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").Value = '3.610.65'
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '66.713.37'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("E19").Value = '  +3.05%  '
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.618'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '78.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '3.746.07'
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.73%  '
$ws.Range("E28").Value = '  +3.18%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").Value = '3.599.50'
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.159'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.35%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.65'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '177.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("E40").Value = '  +0.70%  '
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("E44").Value = '  +9.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.15%  '
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("D51").Value = '2.429.75'
$ws.Range("E51").Value = '  +5.51%  '
